$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Split the combined F1 and F2 cross labels into separate a/b replicate rows.
$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

# Update the active selection to reflect the last edited cell.
$ws.Range("A7").Select()
